# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar el texto de la conversión del día (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$nl = [char]10
$newText = "Conversión del día 💰" + $nl + `
    "✅ Dólar paralelo: 68" + $nl + `
    $nl + `
    "Binance" + $nl + `
    "✅ 1000 Bs = 12.55 = 50288.21 pesos" + $nl + `
    "✅ 50288.21 pesos = 12.52 = 970.46 Bs" + $nl + `
    $nl + `
    "Promedio competencia" + $nl + `
    "✅ Tasa pesos: 20" + $nl + `
    "✅ Tasa Bs: 20" + $nl + `
    "✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: actualizar tasas N10/O10 y N12/O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 79.7
$wsTasas.Range("O10").Value = 4007.97
$wsTasas.Range("N12").Value = 4017
$wsTasas.Range("O12").Value = 77.52
